$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "41.966.50"
$ws.Range("E2").Value = "  -1.17%  "

# Row 3
$ws.Range("D3").Value = "2.245.12"
$ws.Range("E3").Value = "  -1.72%  "

# Row 4
$ws.Range("E4").Value = "  +0.13%  "

# Row 5
$ws.Range("D5").Value = "'305.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.63%  "

# Row 6
$ws.Range("D6").Value = "'96.31"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.74%  "

# Row 7
$ws.Range("D7").Value = "'0.523"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.57%  "

# Row 8
$ws.Range("E8").Value = "  +0.05%  "

# Row 9
$ws.Range("D9").Value = "'0.486"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.80%  "

# Row 10
$ws.Range("D10").Value = "'34.67"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.87%  "

# Row 11
$ws.Range("D11").Value = "'0.0808"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.20%  "

# Row 12
$ws.Range("E12").Value = "  +0.97%  "

# Row 13
$ws.Range("D13").Value = "'6.75"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.49%  "

# Row 14
$ws.Range("D14").Value = "2.600.43"
$ws.Range("E14").Value = "  -1.41%  "

# Row 15
$ws.Range("D15").Value = "'14.35"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.90%  "

# Row 16
$ws.Range("D16").Value = "2.237.91"
$ws.Range("E16").Value = "  -2.61%  "

# Row 17
$ws.Range("D17").Value = "'0.776"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.81%  "

# Row 18
$ws.Range("D18").Value = "41.918.27"
$ws.Range("E18").Value = "  -1.07%  "

# Row 19
$ws.Range("D19").Value = "'12.11"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.75%  "

# Row 20
$ws.Range("D20").Value = "0.0₃0898"
$ws.Range("E20").Value = "  -1.68%  "

# Row 21
$ws.Range("D21").Value = "'5.91"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.36%  "

# Row 22
$ws.Range("D22").Value = "'67.08"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.14%  "

# Row 23
$ws.Range("D23").Value = "'234.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.68%  "

# Row 24
$ws.Range("E24").Value = "  -1.98%  "

# Row 25
$ws.Range("B25").Value = "ImmutableX"
$ws.Range("C25").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D25").Value = "'1.93"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.41%  "

# Row 26
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.01%  "

# Row 27
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").Value = "'23.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.18%  "

# Row 28
$ws.Range("B28").Value = "InjectiveProtocol"
$ws.Range("C28").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D28").Value = "'37.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.49%  "

# Row 29
$ws.Range("D29").Value = "'2.11"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.64%  "

# Row 30
$ws.Range("D30").Value = "'9.45"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.22%  "

# Row 31
$ws.Range("D31").Value = "'167.39"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.98%  "

# Row 32
$ws.Range("D32").Value = "'1.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.28%  "

# Row 33
$ws.Range("D33").Value = "'5.14"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.40%  "

# Row 34
$ws.Range("D34").Value = "'3.07"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.58%  "

# Row 35
$ws.Range("D35").Value = "'17.37"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.47%  "

# Row 36
$ws.Range("D36").Value = "'0.0717"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.64%  "

# Row 37
$ws.Range("E37").Value = "  +0.20%  "

# Row 38
$ws.Range("D38").Value = "'0.114"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.59%  "

# Row 39
$ws.Range("D39").Value = "'0.102"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -3.34%  "

# Row 40
$ws.Range("D40").Value = "'1.78"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.74%  "

# Row 41
$ws.Range("D41").Value = "'4.03"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.15%  "

# Row 42
$ws.Range("D42").Value = "1.940.50"
$ws.Range("E42").Value = "  -3.06%  "

# Row 43
$ws.Range("D43").Value = "'0.0280"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.25%  "

# Row 44
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").Value = "'2.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -10.64%  "

# Row 45
$ws.Range("B45").Value = "EnergySwap"
$ws.Range("C45").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D45").Value = "'18.41"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.88%  "

# Row 46
$ws.Range("D46").Value = "'2.88"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.72%  "

# Row 47
$ws.Range("D47").Value = "'9.61"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.02%  "

# Row 48
$ws.Range("D48").Value = "'53.65"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.03%  "

# Row 49
$ws.Range("D49").Value = "2.467.68"
$ws.Range("E49").Value = "  -1.51%  "

# Row 50
$ws.Range("D50").Value = "'70.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.80%  "

# Row 51
$ws.Range("D51").Value = "'90.80"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.72%  "
